$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: drop the obsolete remark "actually venous pH for now" (F3) ---
$ws.Range("F3").ClearContents()

# --- Row 29: was an empty merged placeholder row; replace with a real
#     "thrombocytes" row (lab variable added in the 06/2022 extraction) ---
$ws.Range("A29:E29").UnMerge()

$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 1000
$ws.Cells.Item(29, 5).NumberFormat = "#,##0"

$ws.Cells.Item(29, 1).Value = "thrombocytes"
$ws.Cells.Item(29, 1).Style = "Normal"
$ws.Cells.Item(29, 2).Value = "thrombocytes"
$ws.Cells.Item(29, 2).Style = "Normal"

$ws.Cells.Item(29, 3).Value = "G/l"
$ws.Cells.Item(29, 3).NumberFormat = "0"
$ws.Cells.Item(29, 3).VerticalAlignment = -4107

# --- Selection cursor, as last saved ---
$ws.Range("B32").Select() | Out-Null
